$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 7; existing rows 7-27 shift down to 8-28.
$ws.Rows("7").Insert()

# Populate the new row 7 with the weekly record (same market/product context
# as the row that used to be row 7, but a new date/volume/price observation).
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "Macroferia Regional de Talca"
$ws.Range("C7").Value = "Maule"
$ws.Range("D7").Value = 44519
$ws.Range("E7").Value = 7
$ws.Range("F7").Value = "Fruta"
$ws.Range("G7").Value = 100101
$ws.Range("H7").Value = "Berries"
$ws.Range("I7").Value = 100101001
$ws.Range("J7").Value = "Arándano (blue)"
$ws.Range("K7").Value = "Sin especificar"
$ws.Range("L7").Value = "Primera"
$ws.Range("M7").Value = 180
$ws.Range("N7").Value = 4000
$ws.Range("O7").Value = 4000
$ws.Range("P7").Value = 4000
$ws.Range("Q7").Value = "$/bandeja 2 kilos"
$ws.Range("R7").Value = "Provincia de Linares"
$ws.Range("S7").Value = 2000
$ws.Range("T7").Value = 2
